# Auto-generated edit script applying the cryptos price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '61.995.36'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -2.32%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.579.81'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -4.04%  '
$ws.Range('E4').Value = '  +0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '549.05'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.01%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '154.69'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -2.33%  '
$ws.Range('E7').Value = '  +0.04%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.594'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +1.92%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.104'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('E10').Value = '  -1.30%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '5.57'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +3.78%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.365'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.92%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '3.037.11'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -3.97%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '25.60'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -2.90%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '61.944.52'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -2.16%  '
$ws.Range('E16').Value = '  -0.31%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '2.582.41'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -4.02%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '11.67'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -3.02%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '4.55'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.49%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '338.05'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.32%  '
$ws.Range('E21').Value = '  -4.41%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.39%  '
$ws.Range('E23').Value = '  -2.76%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '63.74'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -0.27%  '
$ws.Range('E25').Value = '  -0.74%  '
$ws.Range('E26').Value = '  -0.34%  '
$ws.Range('E27').Value = '  +0.05%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '7.30'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +4.20%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '0.0₃0839'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('E30').Value = '  +2.04%  '
$ws.Range('E31').Value = '  -3.19%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '162.94'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -1.41%  '
$ws.Range('E33').Value = '  +2.09%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('E35').Value = '  -1.67%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.42'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -0.57%  '
$ws.Range('E37').Value = '  +1.07%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '329.22'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -3.20%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '6.04'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.30%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.909'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -4.03%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '37.65'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -1.22%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '21.00'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +0.98%  '
$ws.Range('E44').Value = '  -0.01%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.608'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -1.95%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '10.96'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.86%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '2.113.43'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +0.91%  '
$ws.Range('E48').Value = '  -2.68%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '19.58'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -3.46%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0965'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.71%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0238'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -1.03%  '
